$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Simple value fixes in column E
$ws.Range("E12").Value = 1361974148
$ws.Range("E14").Value = 537985.22
$ws.Range("E15").Value = 30000000
$ws.Range("E16").Value = -60473972.810000002

# E18 becomes a formula summing E12:E17 (was a literal value before)
$ws.Range("E18").Formula = "=SUM(E12:E17)"

$ws.Range("E19").Value = -384700000

# E21 becomes a formula summing E18:E20 (was a literal value before)
$ws.Range("E21").Formula = "=SUM(E18:E20)"

$ws.Range("E22").Value = -20015625

# E23 formula already exists (=SUM(E21:E22)); recalculation will update the value.

# E25 formula already exists (=E23/E26); fix E26 so it stops erroring (#VALUE!)
$ws.Range("E26").Value = 1018613404

$wb.Save()
